# Add a new "calldata" worksheet after the existing "contacts" sheet and
# populate it with the call/task/case/note sample rows, matching the
# target workbook produced by the commit "i have added call ctest class
# and extend report and log 4j".

$wb = $excel.ActiveWorkbook

# Insert the new sheet directly after "contacts" so sheet order becomes
# contacts, calldata (Worksheets.Add defaults to inserting *before* the
# active sheet, so we explicitly pass an After: argument).
$contacts = $wb.Worksheets.Item("contacts")
$calldata = $wb.Worksheets.Add([Type]::Missing, $contacts)
$calldata.Name = "calldata"

# Header row.
$calldata.Cells.Item(1, 1).Value = "deal"
$calldata.Cells.Item(1, 2).Value = "task"
$calldata.Cells.Item(1, 3).Value = "case"
$calldata.Cells.Item(1, 4).Value = "note"

# Row 2.
$calldata.Cells.Item(2, 1).Value = "a"
$calldata.Cells.Item(2, 2).Value = "b"
$calldata.Cells.Item(2, 3).Value = "c"
$calldata.Cells.Item(2, 4).Value = "aaa"

# Row 3.
$calldata.Cells.Item(3, 1).Value = "e"
$calldata.Cells.Item(3, 2).Value = "f"
$calldata.Cells.Item(3, 3).Value = "g"
$calldata.Cells.Item(3, 4).Value = "bbb"

# Row 4 - note column written before the rest of the row to mirror the
# shared-string insertion order of the authored workbook.
$calldata.Cells.Item(4, 4).Value = "xxx"
$calldata.Cells.Item(4, 1).Value = "h"
$calldata.Cells.Item(4, 2).Value = "i"
$calldata.Cells.Item(4, 3).Value = "j"

# Match the final selection/active cell left on the calldata sheet.
$calldata.Range("C4").Select()
